# Extend the sanitation-access table with two more years (2021, 2022) in
# columns N and O, mirroring the existing D:M layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build each distinct new cell format once in a scratch area, well outside
# the table, then Copy / PasteSpecial(xlPasteFormats) it onto the real
# target range. Doing it this way (one atomic "apply finished format"
# step per target range) avoids leaving behind orphaned intermediate
# cellXfs entries that a step-by-step Font/NumberFormat assignment over a
# multi-column range would otherwise generate.
$scratchBold    = $ws.Range("Z101")
$scratchRegular = $ws.Range("Z102")
$scratchBorder  = $ws.Range("Z103")

$scratchBold.NumberFormat = "0.0"
$scratchBold.Font.Name = "Times New Roman"
$scratchBold.Font.Size = 10
$scratchBold.Font.Bold = $true

$scratchRegular.NumberFormat = "0.0"
$scratchRegular.Font.Name = "Times New Roman"
$scratchRegular.Font.Size = 10
$scratchRegular.Font.Bold = $false

# For the bordered variant, start from a cell that already carries the
# exact "medium bottom border" (M14) via Copy/PasteSpecial rather than
# building the border up through the Borders collection property-by-
# property -- the latter registers throwaway intermediate border
# definitions (e.g. a transient "thin" weight) that never get cleaned
# back out of the workbook's border table.
$ws.Range("M14").Copy() | Out-Null
$scratchBorder.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$scratchBorder.Font.Size = 10

# --- Row 3 (thin separator row above the header) ------------------------
# N3/O3 just pick up the same "thick bottom border" formatting as K3:M3.
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3:O3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 4 (year header row) ---------------------------------------------
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4:O4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("N4").Value2 = 2021
$ws.Range("O4").Value2 = 2022

# --- Row 5 (Kyrgyz Republic totals row, bold) ----------------------------
$scratchBold.Copy() | Out-Null
$ws.Range("N5:O5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("N5").Value2 = 40.007977647471066
$ws.Range("O5").Value2 = 42.620582506455563

# --- Rows 6-13 (oblast data rows, regular weight) ------------------------
$scratchRegular.Copy() | Out-Null
$ws.Range("N6:O13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("N6").Value2 = 5.7072514621689896
$ws.Range("O6").Value2 = 8.1443914479075037
$ws.Range("N7").Value2 = 8.9893229854028949
$ws.Range("O7").Value2 = 10.715961386284755
$ws.Range("N8").Value2 = 66.307512472824584
$ws.Range("O8").Value2 = 81.977461999426666
$ws.Range("N9").Value2 = 23.475213049310256
$ws.Range("O9").Value2 = 29.828871240443185
$ws.Range("N10").Value2 = 9.8045372040896162
$ws.Range("O10").Value2 = 9.7218425128664112
$ws.Range("N11").Value2 = 9.3737779268960448
$ws.Range("O11").Value2 = 8.6167819403064012
$ws.Range("N12").Value2 = 70.457032471318783
$ws.Range("O12").Value2 = 69.915337594090886
$ws.Range("N13").Value2 = 98.411252120183207
$ws.Range("O13").Value2 = 99.08571752721997

# --- Row 14 (Osh oblast row, bottom border) -------------------------------
$scratchBorder.Copy() | Out-Null
$ws.Range("N14:O14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("N14").Value2 = 63.900563564170795
$ws.Range("O14").Value2 = 64.805252627098838

# --- Clean up the scratch cells so they don't linger in the sheet --------
$ws.Range("Z101:Z103").Clear() | Out-Null

# --- Sheet-level bookkeeping ----------------------------------------------
$ws.Range("P8").Select() | Out-Null
